# "map updates changed the order of almost all the maps in the game
#  updated memory requirements of maps"
#
# The "Maps" table on Sheet1 (header row 121, data rows formerly 122-136)
# got reshuffled:
#   - the old row for map #1 ("ok") is replaced with a brand new map
#     ("short and sweet", 56 bytes) and relocated into the middle of the
#     table (row 131)
#   - maps #6,#7,#8,#9,#10 (old rows 127-131) were re-sorted by their
#     byte/memory requirement (column B) into rows 126-130
#   - maps #2,#3,#4,#5 (old rows 123-126) were moved down to the bottom
#     of the table (rows 137-140)
#   - maps #11-#15 stay where they were (rows 132-136)
#   - map #9's "jump-move forward..." sub-info (columns H/I) moved to
#     map #8, and map #12 picked up a "can't be beaten" note

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old rows 122-125 (maps #1-#4) are fully vacated - their old
# content either moved elsewhere or was replaced, and no map table rows
# remain at 122-125 afterwards.
$ws.Range("A122:L125").Clear()

# --- rows 126-130: maps #7, #6, #10, #8, #9 sorted ascending by column B ---
$ws.Cells.Item(126,1).Value = 7
$ws.Cells.Item(126,2).Value = 64
$ws.Cells.Item(126,6).Value = "very easy, should be an earlier level"

$ws.Cells.Item(127,1).Value = 6
$ws.Cells.Item(127,2).Value = 78
$ws.Cells.Item(127,6).Value = "very easy, should be an earlier level"

$ws.Cells.Item(128,1).Value = 10
$ws.Cells.Item(128,2).Value = 154
$ws.Cells.Item(128,6).Value = "easy - very straightforward"

$ws.Cells.Item(129,1).Value = 8
$ws.Cells.Item(129,2).Value = 156
$ws.Cells.Item(129,6).Value = "fairly easy - fun water map"
$ws.Cells.Item(129,8).Value = "jump-move forward and jump-move forward 3 are EASILY mass repeated on this map, will test with subs"
$ws.Cells.Item(129,9).Value = "sub1 = jump/move forward, sub2 = jump"

$ws.Cells.Item(130,1).Value = 9
$ws.Cells.Item(130,2).Value = 306
$ws.Cells.Item(130,3).Value = 164
$ws.Cells.Item(130,6).Value = "fun - lots of jumping"
$ws.Range("H130:I130").Clear()

# --- row 131: map #1 completely redone (new content, new memory size) ---
$ws.Cells.Item(131,1).Value = 1
$ws.Cells.Item(131,2).Value = 56
$ws.Cells.Item(131,6).Value = "short and sweet"

# --- row 132: map #11, unchanged ---
$ws.Cells.Item(132,1).Value = 11
$ws.Cells.Item(132,2).Value = 128
$ws.Cells.Item(132,6).Value = "not very complex, fairly easy as well - lots of random extra stuff not related to finishing the map"

# --- row 133: map #12 now annotated "can't be beaten" ---
$ws.Cells.Item(133,1).Value = 12
$ws.Cells.Item(133,6).Value = "can't be beaten"

# --- row 134: map #13, unchanged ---
$ws.Cells.Item(134,1).Value = 13
$ws.Cells.Item(134,2).Value = 182
$ws.Cells.Item(134,6).Value = "could probably lose the reprogram square, additionally, not sure if intended, but the last two switches can be skipped entirely, may want to disable jump on this map, or make some reason for the switches"

# --- row 135: map #14, unchanged ---
$ws.Cells.Item(135,1).Value = 14
$ws.Cells.Item(135,2).Value = 318
$ws.Cells.Item(135,6).Value = "very linear, interesting figuring out what does what"
$ws.Cells.Item(135,12).Value = "needs edge squares removed maybe"

# --- row 136: map #15, unchanged ---
$ws.Cells.Item(136,1).Value = 15
$ws.Cells.Item(136,6).Value = "can't be beaten"

# --- rows 137-140: maps #2,#3,#4,#5 moved to the bottom of the table ---
$ws.Cells.Item(137,1).Value = 2
$ws.Cells.Item(137,2).Value = 156
$ws.Cells.Item(137,6).Value = "annoying as hell"

$ws.Cells.Item(138,1).Value = 3
$ws.Cells.Item(138,2).Value = 110
$ws.Cells.Item(138,4).Value = "62 with reprogram"
$ws.Cells.Item(138,6).Value = "annoying as hell"

$ws.Cells.Item(139,1).Value = 4
$ws.Cells.Item(139,2).Value = 88
$ws.Cells.Item(139,6).Value = "annoying as hell"
$ws.Cells.Item(139,11).Value = "used left switch"

$ws.Cells.Item(140,1).Value = 5
$ws.Cells.Item(140,2).Value = 92
$ws.Cells.Item(140,6).Value = "annoying as hell"
$ws.Cells.Item(140,11).Value = "fixed bad switch"

# The sort that produced rows 126-130 was a "Data > Sort" on that block,
# keyed on the byte/memory column.
$ws.Range("A126:G130").Sort($ws.Range("B126:B130"))

# Reflect the resulting selection/viewport like the author left it.
$ws.Range("F132").Select()
